$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: ht=45, B/C/D/E filled ---
$ws.Cells.Item(29, 2).Value = 'DiscordiaAgency_Demo_2017_09_23.exe'
$ws.Range("B28").Copy()
$ws.Cells.Item(29, 2).PasteSpecial(-4122)

$ws.Cells.Item(29, 3).Value = "Entwicklung"
$ws.Range("C28").Copy()
$ws.Cells.Item(29, 3).PasteSpecial(-4122)

$ws.Cells.Item(29, 4).Value = "Anna Franziska"
$ws.Range("D28").Copy()
$ws.Cells.Item(29, 4).PasteSpecial(-4122)

$ws.Cells.Item(29, 5).Value = 'Intro-Text im Startmenü; Hilfe-Seite mit Steuerung; Wachen bleiben jetzt mit etwas Abstand zum aufgetroffenen Objekt'
$ws.Range("E28").Copy()
$ws.Cells.Item(29, 5).PasteSpecial(-4122)

$ws.Rows.Item(29).RowHeight = 45

# --- Row 30: ht=210, A/B/C/D/F filled ---
$ws.Cells.Item(30, 1).Value = 43001
$ws.Range("A26").Copy()
$ws.Cells.Item(30, 1).PasteSpecial(-4122)

$ws.Cells.Item(30, 2).Value = "DiscordiaAgency_Demo_2017_09_22-3.exe"
$ws.Range("B28").Copy()
$ws.Cells.Item(30, 2).PasteSpecial(-4122)

$ws.Cells.Item(30, 3).Value = "Spielen"
$ws.Range("C28").Copy()
$ws.Cells.Item(30, 3).PasteSpecial(-4122)

$ws.Cells.Item(30, 4).Value = "Urban"
$ws.Range("D28").Copy()
$ws.Cells.Item(30, 4).PasteSpecial(-4122)

$ws.Cells.Item(30, 6).Value = 'dragged Guard still collides with other dead Guards; speed-run-achievements and level timer; bullets still rather easy to avoid; Level 7: very easy, more guards, or 2 guards on same path; guards they should say "hello" or "yo dawg" when meeting each other; try Deus Ex & Commandos 2 as inspiration; "why are you always going towards the exit?" -> Target Sprite needs to be changed, probably, or use them as exits and chain levels, so you need to go through several rooms, before actually seeing an assassination target that you have to specifically kill somehow; maybe have "strong" guards that cannot be disabled; or cameras, that can''t shoot you, but alert some guards to you; [12:50] [Group] [Anuschka]: btw. are the level intro fonts still too large for your screen?
[12:50] [Group] [*Pel]: hmm it''s better at least
[12:51] [Group] [*Pel]: and when guards yell "argh" it should alert nearby guards in a small radius :p
[12:51] [Group] [*Pel]: but that would affec t your level design
[12:51] [Group] [*Pel]: they do yell quite loud...; [12:59] [Group] [*Pel]: and a narrator like in stanley''s parable
[12:59] [Group] [Anuschka]: lol
[13:00] [Group] [*Pel]: "[player] hit the poor guard John on the back and is now struggling with his moral conscience"
[13:00] [Group] [*Pel]: "does that make him a bad guy"'
$ws.Range("F24").Copy()
$ws.Cells.Item(30, 6).PasteSpecial(-4122)

$ws.Rows.Item(30).RowHeight = 210

# --- Row 31: ht=105, B/C/D/E filled ---
$ws.Cells.Item(31, 2).Value = 'DiscordiaAgency_Demo_2017_09_23-2.exe'
$ws.Range("B28").Copy()
$ws.Cells.Item(31, 2).PasteSpecial(-4122)

$ws.Cells.Item(31, 3).Value = "Entwicklung"
$ws.Range("C28").Copy()
$ws.Cells.Item(31, 3).PasteSpecial(-4122)

$ws.Cells.Item(31, 4).Value = "Anna Franziska"
$ws.Range("D28").Copy()
$ws.Cells.Item(31, 4).PasteSpecial(-4122)

$ws.Cells.Item(31, 5).Value = 'Tote Wachen können nicht mehr miteinander kollidieren; Bug: angezeigter Radius der Geräuschquelle stimmt nicht mit Hörweite überein; Bug: falsche Ausgangsposition für Rufradius der Wachen; Wachen rotieren leicht nach links und rechts, nachdem sie Geräusch gefolgt sind; stationäre Wachen rotieren leicht nach links und rechts'
$ws.Range("E28").Copy()
$ws.Cells.Item(31, 5).PasteSpecial(-4122)

$ws.Rows.Item(31).RowHeight = 105

# --- sheet view / dimension updates ---
$null = $ws.Range("E32").Select()

Write-Host "done"
